# fix(publipostage): Try to solve Excel emoji problem
#
# The workbook used four emoji icons as "statut" markers (column A):
#   blue book   (U+1F4D8) -> warning sign "⚠️" (U+26A0 U+FE0F)
#   red book    (U+1F4D5) -> "-3"
#   orange book (U+1F4D9) -> "+3"
#   green book  (U+1F4D7) -> check mark "✅" (U+2705)
#
# This mirrors the shared-string table edit: every cell that used to
# contain one of the old emoji now contains the corresponding new
# text/emoji; nothing else in the workbook (other shared strings,
# styles, other cells) is touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldBlue   = "📘"
$oldRed    = "📕"
$oldOrange = "📙"
$oldGreen  = "📗"

$newWarning = "⚠️"
$newMinus3  = "-3"
$newPlus3   = "+3"
$newCheck   = "✅"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Cells whose replacement text looks like a plain number ("-3" / "+3")
# must be entered in a way that keeps them as *text*, otherwise Excel's
# automatic type inference would silently store them as numeric values
# instead of strings (and would also implicitly change the cell's number
# format/style). Using a formula that yields a text result, then
# converting that formula to its static value via Copy / Paste-Special-
# Values, achieves this without altering the cell's number format.
for ($r = 1; $r -le $lastRow; $r++) {
  for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    $val = $cell.Value()

    # Only text cells are candidates. Skip booleans/numbers/blanks
    # outright: comparing them with "-eq" against a (non-empty) string
    # literal is unsafe in PowerShell, because "-eq" coerces the right
    # operand to the type of the left operand, so e.g. $true -eq "📘"
    # would otherwise (incorrectly) evaluate to $true.
    if ($val -isnot [string]) {
      continue
    }

    if ($val -ceq $oldBlue) {
      $cell.Value = $newWarning
    } elseif ($val -ceq $oldGreen) {
      $cell.Value = $newCheck
    } elseif ($val -ceq $oldRed) {
      $cell.Formula = '="-3"'
      $cell.Copy()
      $cell.PasteSpecial(-4163)   # xlPasteValues
    } elseif ($val -ceq $oldOrange) {
      $cell.Formula = '="+3"'
      $cell.Copy()
      $cell.PasteSpecial(-4163)   # xlPasteValues
    }
  }
}

$excel.CutCopyMode = $false
